$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "2026-02-15 03:55:50"
$ws.Range("C10").Value = "LOUISE DAJEU"
$ws.Range("D10").Value = 7518

# Row 11
$ws.Range("A11").Value = "2026-02-15 03:00:00"
$ws.Range("C11").Value = "ROSE DEUMENI"
$ws.Range("D11").Value = 14870

# Row 13
$ws.Range("A13").Value = "2026-02-15 06:12:30"
$ws.Range("C13").Value = "Clarisse Ngenue Wankah"
$ws.Range("D13").Value = 12048

# Row 15
$ws.Range("A15").Value = "2026-02-15 02:07:22"
$ws.Range("D15").Value = 35306

# Row 17
$ws.Range("A17").Value = "2026-02-15 06:04:15"
$ws.Range("C17").Value = "ETIENNE JUSTIN JIOFACK"
$ws.Range("D17").Value = 2717

# Row 18
$ws.Range("A18").Value = "2026-02-15 05:17:39"
$ws.Range("C18").Value = "CLARISSE MAKOLO"
$ws.Range("D18").Value = 7758

# Row 19
$ws.Range("A19").Value = "2026-02-15 06:08:06"
$ws.Range("C19").Value = "LANDRY MANFOUO"
$ws.Range("D19").Value = 3855
